# "improved P2 based on peer review feedback"
# The weekly "Period" header row (F8:AM8) used to show absolute calendar
# dates; it is changed here to show the ISO week number of each period
# instead (matching the "WEEKS" column header already in F7), and the
# cell number format is switched from the custom date format to a plain
# integer format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project")

# New period (week-number) values for F8:AM8, replacing the old calendar
# date serials.
$periodCols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM")
$periodWeeks = @(46,47,48,49,50,51,52,53,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26)

for ($i = 0; $i -lt $periodCols.Length; $i++) {
    $ws.Range($periodCols[$i] + "8").Value = $periodWeeks[$i]
}

# Switch the period row from the custom date format to a plain integer
# format now that it holds week numbers rather than dates.
$ws.Range("F8:AM8").NumberFormat = "0"

# Leave the cursor where the author left it after the edit.
$ws.Range("C27").Select()
